# Apply the edit described by the commit:
#  - rename the sheet to "Aiden Markram"
#  - insert a new "matchNo" column at the front
#  - populate the full batting log (6 innings) for Aiden Markram

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Aiden Markram"

# Shift the existing columns (teamName..result) one place to the right
# and make room for the new leading "matchNo" column.
$ws.Range("A1").EntireColumn.Insert()

# Cells whose text looks like a plain number ("13", "162.50", ...) must be
# written with a Text number format first, otherwise Excel/COM would silently
# coerce them into real numbers and drop things like trailing zeros.
function Set-TextValue($ws, $row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$headers = @(
    "matchNo", "teamName", "batterName", "states", "runs", "balls", "fours", "sixes", "sr", "opponentTeamName", "venue", "date", "result"
)
for ($c = 0; $c -lt $headers.Count; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# matchNo, teamName, batterName, states, runs, balls, fours, sixes, sr,
# opponentTeamName, venue, date, result
$data = @(
    @("53rd", "Punjab Kings", "Aiden Markram", "c †Dhoni b Thakur", "13", "8", "0", "1", "162.50", "Chennai Super Kings", "Dubai (DSC)", "October 07", "Punjab Kings won by 6 wickets (with 42 balls remaining)"),
    @("37th", "Punjab Kings", "Aiden Markram", "c Pandey b Abdul Samad", "27", "32", "2", "0", "84.37", "Sunrisers Hyderabad", "Sharjah", "September 25", "Punjab Kings won by 5 runs"),
    @("42nd", "Punjab Kings", "Aiden Markram", "b Chahar", "42", "29", "6", "0", "144.82", "Mumbai Indians", "Abu Dhabi", "September 28", "Mumbai won by 6 wickets (with 6 balls remaining)"),
    @("45th", "Punjab Kings", "Aiden Markram", "c Shubman Gill b Narine", "18", "16", "0", "1", "112.50", "Kolkata Knight Riders", "Dubai (DSC)", "October 01", "Punjab Kings won by 5 wickets (with 3 balls remaining)"),
    @("32nd", "Punjab Kings", "Aiden Markram", "", "26", "20", "2", "1", "130.00", "Rajasthan Royals", "Dubai (DSC)", "September 21", "Royals won by 2 runs"),
    @("48th", "Punjab Kings", "Aiden Markram", "c Christian b Garton", "20", "14", "2", "1", "142.85", "Royal Challengers Bangalore", "Sharjah", "October 03", "RCB won by 6 runs")
)

# Columns (1-based) whose values must be forced to Text so Excel keeps them
# exactly as scraped (runs, balls, fours, sixes, sr).
$numericTextCols = @(5, 6, 7, 8, 9)

for ($r = 0; $r -lt $data.Count; $r++) {
    $row = $data[$r]
    $excelRow = $r + 2
    for ($c = 0; $c -lt $row.Count; $c++) {
        $excelCol = $c + 1
        if ($numericTextCols -contains $excelCol) {
            Set-TextValue $ws $excelRow $excelCol $row[$c]
        } else {
            $ws.Cells.Item($excelRow, $excelCol).Value = $row[$c]
        }
    }
}

Write-Host "Updated sheet now spans $($ws.UsedRange.Address())"
